## plano de aula revisado
##
## The four class dates in column A for rows 25-28 (4, 6, 11 and 13 Nov)
## were stored as real Excel date serials (45600/45602/45607/45609). They
## are switched to plain text labels, matching the "d-Mmm" style already
## used by every other date in the column (e.g. "30-Out" in A24). A leading
## backtick is typed so the text is stored verbatim instead of being
## re-parsed into a date value.
##
## Scroll the sheet up one row (view was sitting at A17, now A16) to match
## the refreshed layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$excel.ActiveWindow.ScrollRow = 16

$ws.Range("A25").Value = '`4-Nov'
$ws.Range("A26").Value = '`6-Nov'
$ws.Range("A27").Value = '`11-Nov'
$ws.Range("A28").Value = '`13-Nov'
